$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling "serorreversion" -> "seroreversion" in the second (A8:G12) table headers
$ws.Range("F8").Value = "Tasa de seroreversion"
$ws.Range("G8").Value = "Tasa de seroreversión Rhat"

# Update the active selection to the second table range
$ws.Range("A8:G12").Select()
